$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.548.42'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").Value = '2.292.68'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.76'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.96'
$ws.Range("E6").Value = '  -2.96%  '
$ws.Range("E7").Value = '  -0.79%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -1.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.51'
$ws.Range("E10").Value = '  -3.69%  '
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.55'
$ws.Range("E12").Value = '  -5.69%  '
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.30'
$ws.Range("E14").Value = '  +4.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.74'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '2.648.21'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").Value = '2.275.83'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.795'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").Value = '42.464.01'
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.72'
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("E21").Value = '  -0.75%  '
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.74'
$ws.Range("E23").Value = '  -1.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.93'
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("E25").Value = '  +0.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.46'
$ws.Range("E26").Value = '  -1.73%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.93'
$ws.Range("E28").Value = '  -3.26%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.69'
$ws.Range("E30").Value = '  +1.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.19'
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.19'
$ws.Range("E32").Value = '  +1.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.71'
$ws.Range("E34").Value = '  +5.11%  '
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.01'
$ws.Range("E36").Value = '  +1.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0698'
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("E38").Value = '  -3.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.82'
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0995'
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("E42").Value = '  -1.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.28'
$ws.Range("E43").Value = '  -9.38%  '
$ws.Range("D44").Value = '1.966.27'
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0279'
$ws.Range("E45").Value = '  -0.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.65'
$ws.Range("E46").Value = '  -3.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.71'
$ws.Range("E47").Value = '  -5.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.83'
$ws.Range("E48").Value = '  -1.65%  '
$ws.Range("D49").Value = '2.515.11'
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.69'
$ws.Range("E50").Value = '  -4.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.75'
$ws.Range("E51").Value = '  -2.82%  '
